$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

function Set-TextCell($sheet, $addr, $val) {
    # Force the cell to be stored as text (matching the "numberStoredAsText"
    # convention already used throughout this workbook) rather than as a
    # number, even when the value looks numeric (e.g. "5", "17", "3").
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($addr).Style = "Normal"
}

# New order rows appended to the "Orders" sheet (rows 22-31)
$ws.Range("C22").Value = "7_翠绿洋桔梗_Dark Green Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextCell $ws "F22" "5"

$ws.Range("C23").Value = "13_酒红洋桔梗_Burgundy Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextCell $ws "F23" "5"

$ws.Range("C24").Value = "300_白星_White Gypso_ gypsophila_1kg"
Set-TextCell $ws "F24" "5"

Set-TextCell $ws "A25" "3"
$ws.Range("C25").Value = "597_尤加利叶小叶_undefined_undefined_1bunch"
Set-TextCell $ws "F25" "17"

$ws.Range("C26").Value = "321_雪柳叶_Spiraea  leaves_undefined_1bunch"
Set-TextCell $ws "F26" "10"

$ws.Range("C27").Value = "542_吊米 红_hanging amaranthus`nred_undefined_1bunch"
Set-TextCell $ws "F27" "10"

$ws.Range("C28").Value = "3_波浪白洋桔梗_Wavy White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextCell $ws "F28" "15"

Set-TextCell $ws "A29" "4"
$ws.Range("C29").Value = "321_雪柳叶_Spiraea  leaves_undefined_1bunch"
Set-TextCell $ws "F29" "10"

$ws.Range("C30").Value = "371_红朱蕉_Cordyline red_undefined_1bunch"
Set-TextCell $ws "F30" "5"

$ws.Range("C31").Value = "320_雪柳花_Spiraea flower white_undefined_1bunch"

# Update the running total number string on the "Summary" sheet
Set-TextCell $summary "G2" "0588103102020555851031215655555171010151050"
